# Fruta / hortaliza, semanal
# Weekly refresh: the data rows (2-27) get reshuffled. For each destination
# row, the columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio
# maximo), M (Precio promedio ponderado) and P (Precio $/Kg) are replaced
# with the values that used to live in another ("source") row, while every
# other column stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of destination row -> source row (values are copied FROM source row
# INTO destination row, for columns D, J, K, L, M, P).
$rowMap = @{
    2  = 15
    3  = 20
    4  = 24
    5  = 13
    6  = 17
    7  = 27
    8  = 22
    9  = 5
    10 = 2
    11 = 12
    12 = 7
    13 = 23
    14 = 11
    15 = 19
    16 = 9
    17 = 14
    18 = 16
    19 = 18
    20 = 10
    21 = 25
    22 = 6
    23 = 26
    24 = 21
    25 = 8
    26 = 4
    27 = 3
}

# Columns (by index) that move together with each row in the reshuffle.
# D=4, J=10, K=11, L=12, M=13, P=16
$cols = @(4, 10, 11, 12, 13, 16)

# First, snapshot the current (pre-edit) values for every affected column
# in every row, so that writes to one destination row never clobber data
# that still needs to be read as a source for another row.
$original = @{}
foreach ($row in 2..27) {
    $original[$row] = @{}
    foreach ($col in $cols) {
        $original[$row][$col] = $ws.Cells.Item($row, $col).Value2
    }
}

# Now write the shuffled values into each destination row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $ws.Cells.Item($destRow, $col).Value2 = $original[$srcRow][$col]
    }
}

$wb.Save()
